$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(375, 44449, 1, 4, 175.1313485113835),
    @(376, 44450, 0, 3, 131.3485113835376),
    @(377, 44451, 0, 3, 131.3485113835376),
    @(378, 44452, 0, 2, 87.56567425569177),
    @(379, 44453, 1, 3, 131.3485113835376),
    @(380, 44454, 0, 3, 131.3485113835376),
    @(381, 44455, 0, 2, 87.56567425569177),
    @(382, 44456, 0, 1, 43.78283712784589),
    @(383, 44457, 0, 1, 43.78283712784589),
    @(384, 44458, 0, 1, 43.78283712784589),
    @(385, 44459, 0, 1, 43.78283712784589)
)

# Copy the date-column formatting (bold, bordered, centered, date number format)
# from the last existing date cell (A374) onto the new date cells (A375:A385)
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)  # xlPasteFormats

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}
